$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 84,9
$data[0,0] = 24
$data[0,1] = "Cromossoma 24"
$data[0,2] = 1258.155462688554
$data[0,3] = 14
$data[0,4] = 83
$data[0,5] = 83
$data[0,6] = 166
$data[0,7] = 0.01943566327127971
$data[0,8] = 0.01943566327127971
$data[1,0] = 53
$data[1,1] = "Cromossoma 53"
$data[1,2] = 1263.164564607314
$data[1,3] = 14
$data[1,4] = 82
$data[1,5] = 83
$data[1,6] = 165
$data[1,7] = 0.01931858096241658
$data[1,8] = 0.03875424423369629
$data[2,0] = 13
$data[2,1] = "Cromossoma 13"
$data[2,2] = 1307.413874005071
$data[2,3] = 14
$data[2,4] = 80
$data[2,5] = 83
$data[2,6] = 163
$data[2,7] = 0.01908441634469032
$data[2,8] = 0.05783866057838661
$data[3,0] = 72
$data[3,1] = "Cromossoma 72"
$data[3,2] = 1330.168363200362
$data[3,3] = 14
$data[3,4] = 78
$data[3,5] = 83
$data[3,6] = 161
$data[3,7] = 0.01885025172696406
$data[3,8] = 0.07668891230535066
$data[4,0] = 22
$data[4,1] = "Cromossoma 22"
$data[4,2] = 1344.523447276128
$data[4,3] = 14
$data[4,4] = 75
$data[4,5] = 83
$data[4,6] = 158
$data[4,7] = 0.01849900480037466
$data[4,8] = 0.09518791710572533
$data[5,0] = 57
$data[5,1] = "Cromossoma 57"
$data[5,2] = 1346.924943076636
$data[5,3] = 14
$data[5,4] = 74
$data[5,5] = 83
$data[5,6] = 157
$data[5,7] = 0.01838192249151153
$data[5,8] = 0.1135698395972369
$data[6,0] = 82
$data[6,1] = "Cromossoma 82"
$data[6,2] = 1354.430168242386
$data[6,3] = 14
$data[6,4] = 72
$data[6,5] = 83
$data[6,6] = 155
$data[6,7] = 0.01814775787378527
$data[6,8] = 0.1317175974710221
$data[7,0] = 23
$data[7,1] = "Cromossoma 23"
$data[7,2] = 1362.450290218126
$data[7,3] = 14
$data[7,4] = 71
$data[7,5] = 83
$data[7,6] = 154
$data[7,7] = 0.01803067556492214
$data[7,8] = 0.1497482730359443
$data[8,0] = 56
$data[8,1] = "Cromossoma 56"
$data[8,2] = 1367.452258903949
$data[8,3] = 14
$data[8,4] = 70
$data[8,5] = 83
$data[8,6] = 153
$data[8,7] = 0.01791359325605901
$data[8,8] = 0.1676618662920033
$data[9,0] = 76
$data[9,1] = "Cromossoma 76"
$data[9,2] = 1371.464250490414
$data[9,3] = 14
$data[9,4] = 68
$data[9,5] = 83
$data[9,6] = 151
$data[9,7] = 0.01767942863833275
$data[9,8] = 0.185341294930336
$data[10,0] = 40
$data[10,1] = "Cromossoma 40"
$data[10,2] = 1373.644256429958
$data[10,3] = 14
$data[10,4] = 66
$data[10,5] = 83
$data[10,6] = 149
$data[10,7] = 0.01744526402060649
$data[10,8] = 0.2027865589509425
$data[11,0] = 51
$data[11,1] = "Cromossoma 51"
$data[11,2] = 1378.550347287579
$data[11,3] = 14
$data[11,4] = 64
$data[11,5] = 83
$data[11,6] = 147
$data[11,7] = 0.01721109940288023
$data[11,8] = 0.2199976583538227
$data[12,0] = 58
$data[12,1] = "Cromossoma 58"
$data[12,2] = 1381.931015399535
$data[12,3] = 14
$data[12,4] = 61
$data[12,5] = 83
$data[12,6] = 144
$data[12,7] = 0.01685985247629083
$data[12,8] = 0.2368575108301136
$data[13,0] = 0
$data[13,1] = "Cromossoma 0"
$data[13,2] = 1382.172990947954
$data[13,3] = 14
$data[13,4] = 60
$data[13,5] = 83
$data[13,6] = 143
$data[13,7] = 0.0167427701674277
$data[13,8] = 0.2536002809975413
$data[14,0] = 71
$data[14,1] = "Cromossoma 71"
$data[14,2] = 1220.797969397157
$data[14,3] = 15
$data[14,4] = 84
$data[14,5] = 58
$data[14,6] = 142
$data[14,7] = 0.01662568785856457
$data[14,8] = 0.2702259688561058
$data[15,0] = 63
$data[15,1] = "Cromossoma 63"
$data[15,2] = 1313.465141892968
$data[15,3] = 15
$data[15,4] = 79
$data[15,5] = 58
$data[15,6] = 137
$data[15,7] = 0.01604027631424892
$data[15,8] = 0.2862662451703548
$data[16,0] = 59
$data[16,1] = "Cromossoma 59"
$data[16,2] = 1338.042257335269
$data[16,3] = 15
$data[16,4] = 77
$data[16,5] = 58
$data[16,6] = 135
$data[16,7] = 0.01580611169652266
$data[16,8] = 0.3020723568668774
$data[17,0] = 28
$data[17,1] = "Cromossoma 28"
$data[17,2] = 1343.663377110486
$data[17,3] = 15
$data[17,4] = 76
$data[17,5] = 58
$data[17,6] = 134
$data[17,7] = 0.01568902938765953
$data[17,8] = 0.3177613862545369
$data[18,0] = 60
$data[18,1] = "Cromossoma 60"
$data[18,2] = 1409.717034799425
$data[18,3] = 14
$data[18,4] = 46
$data[18,5] = 83
$data[18,6] = 129
$data[18,7] = 0.01510361784334387
$data[18,8] = 0.3328650040978808
$data[19,0] = 79
$data[19,1] = "Cromossoma 79"
$data[19,2] = 1411.392727747084
$data[19,3] = 13
$data[19,4] = 43
$data[19,5] = 84
$data[19,6] = 127
$data[19,7] = 0.01486945322561761
$data[19,8] = 0.3477344573234984
$data[20,0] = 74
$data[20,1] = "Cromossoma 74"
$data[20,2] = 1369.494088230404
$data[20,3] = 15
$data[20,4] = 69
$data[20,5] = 58
$data[20,6] = 127
$data[20,7] = 0.01486945322561761
$data[20,8] = 0.362603910549116
$data[21,0] = 9
$data[21,1] = "Cromossoma 9"
$data[21,2] = 1372.128282460611
$data[21,3] = 15
$data[21,4] = 67
$data[21,5] = 58
$data[21,6] = 125
$data[21,7] = 0.01463528860789135
$data[21,8] = 0.3772391991570074
$data[22,0] = 35
$data[22,1] = "Cromossoma 35"
$data[22,2] = 1374.685505863701
$data[22,3] = 15
$data[22,4] = 65
$data[22,5] = 58
$data[22,6] = 123
$data[22,7] = 0.01440112399016509
$data[22,8] = 0.3916403231471725
$data[23,0] = 64
$data[23,1] = "Cromossoma 64"
$data[23,2] = 1428.246702605157
$data[23,3] = 14
$data[23,4] = 37
$data[23,5] = 83
$data[23,6] = 120
$data[23,7] = 0.01404987706357569
$data[23,8] = 0.4056902002107481
$data[24,0] = 61
$data[24,1] = "Cromossoma 61"
$data[24,2] = 1380.108853356158
$data[24,3] = 15
$data[24,4] = 62
$data[24,5] = 58
$data[24,6] = 120
$data[24,7] = 0.01404987706357569
$data[24,8] = 0.4197400772743238
$data[25,0] = 78
$data[25,1] = "Cromossoma 78"
$data[25,2] = 1429.932469014709
$data[25,3] = 14
$data[25,4] = 36
$data[25,5] = 83
$data[25,6] = 119
$data[25,7] = 0.01393279475471256
$data[25,8] = 0.4336728720290364
$data[26,0] = 36
$data[26,1] = "Cromossoma 36"
$data[26,2] = 1384.668357435472
$data[26,3] = 15
$data[26,4] = 59
$data[26,5] = 58
$data[26,6] = 117
$data[26,7] = 0.0136986301369863
$data[26,8] = 0.4473715021660227
$data[27,0] = 29
$data[27,1] = "Cromossoma 29"
$data[27,2] = 1385.699726011023
$data[27,3] = 15
$data[27,4] = 58
$data[27,5] = 58
$data[27,6] = 116
$data[27,7] = 0.01358154782812317
$data[27,8] = 0.4609530499941458
$data[28,0] = 70
$data[28,1] = "Cromossoma 70"
$data[28,2] = 1388.209088885934
$data[28,3] = 15
$data[28,4] = 57
$data[28,5] = 58
$data[28,6] = 115
$data[28,7] = 0.01346446551926004
$data[28,8] = 0.4744175155134059
$data[29,0] = 18
$data[29,1] = "Cromossoma 18"
$data[29,2] = 1390.011352515372
$data[29,3] = 15
$data[29,4] = 56
$data[29,5] = 58
$data[29,6] = 114
$data[29,7] = 0.01334738321039691
$data[29,8] = 0.4877648987238028
$data[30,0] = 41
$data[30,1] = "Cromossoma 41"
$data[30,2] = 1390.970913865467
$data[30,3] = 15
$data[30,4] = 55
$data[30,5] = 58
$data[30,6] = 113
$data[30,7] = 0.01323030090153378
$data[30,8] = 0.5009951996253366
$data[31,0] = 42
$data[31,1] = "Cromossoma 42"
$data[31,2] = 1397.822130630478
$data[31,3] = 15
$data[31,4] = 54
$data[31,5] = 58
$data[31,6] = 112
$data[31,7] = 0.01311321859267065
$data[31,8] = 0.5141084182180072
$data[32,0] = 44
$data[32,1] = "Cromossoma 44"
$data[32,2] = 1399.69908564188
$data[32,3] = 15
$data[32,4] = 53
$data[32,5] = 58
$data[32,6] = 111
$data[32,7] = 0.01299613628380752
$data[32,8] = 0.5271045545018147
$data[33,0] = 11
$data[33,1] = "Cromossoma 11"
$data[33,2] = 1401.026717147551
$data[33,3] = 15
$data[33,4] = 52
$data[33,5] = 58
$data[33,6] = 110
$data[33,7] = 0.01287905397494439
$data[33,8] = 0.5399836084767591
$data[34,0] = 20
$data[34,1] = "Cromossoma 20"
$data[34,2] = 1451.610472449456
$data[34,3] = 14
$data[34,4] = 27
$data[34,5] = 83
$data[34,6] = 110
$data[34,7] = 0.01287905397494439
$data[34,8] = 0.5528626624517035
$data[35,0] = 33
$data[35,1] = "Cromossoma 33"
$data[35,2] = 1402.29092511053
$data[35,3] = 15
$data[35,4] = 51
$data[35,5] = 58
$data[35,6] = 109
$data[35,7] = 0.01276197166608125
$data[35,8] = 0.5656246341177847
$data[36,0] = 46
$data[36,1] = "Cromossoma 46"
$data[36,2] = 1451.976195572162
$data[36,3] = 14
$data[36,4] = 26
$data[36,5] = 83
$data[36,6] = 109
$data[36,7] = 0.01276197166608125
$data[36,8] = 0.578386605783866
$data[37,0] = 55
$data[37,1] = "Cromossoma 55"
$data[37,2] = 1405.501169108064
$data[37,3] = 15
$data[37,4] = 50
$data[37,5] = 58
$data[37,6] = 108
$data[37,7] = 0.01264488935721812
$data[37,8] = 0.5910314951410841
$data[38,0] = 8
$data[38,1] = "Cromossoma 8"
$data[38,2] = 1406.595730343208
$data[38,3] = 15
$data[38,4] = 49
$data[38,5] = 58
$data[38,6] = 107
$data[38,7] = 0.01252780704835499
$data[38,8] = 0.6035593021894391
$data[39,0] = 49
$data[39,1] = "Cromossoma 49"
$data[39,2] = 1406.787382723093
$data[39,3] = 15
$data[39,4] = 48
$data[39,5] = 58
$data[39,6] = 106
$data[39,7] = 0.01241072473949186
$data[39,8] = 0.615970026928931
$data[40,0] = 6
$data[40,1] = "Cromossoma 6"
$data[40,2] = 1408.172546050421
$data[40,3] = 15
$data[40,4] = 47
$data[40,5] = 58
$data[40,6] = 105
$data[40,7] = 0.01229364243062873
$data[40,8] = 0.6282636693595597
$data[41,0] = 77
$data[41,1] = "Cromossoma 77"
$data[41,2] = 1469.633976306856
$data[41,3] = 14
$data[41,4] = 21
$data[41,5] = 83
$data[41,6] = 104
$data[41,7] = 0.0121765601217656
$data[41,8] = 0.6404402294813253
$data[42,0] = 2
$data[42,1] = "Cromossoma 2"
$data[42,2] = 1475.217484817921
$data[42,3] = 14
$data[42,4] = 20
$data[42,5] = 83
$data[42,6] = 103
$data[42,7] = 0.01205947781290247
$data[42,8] = 0.6524997072942278
$data[43,0] = 65
$data[43,1] = "Cromossoma 65"
$data[43,2] = 1410.030757386952
$data[43,3] = 15
$data[43,4] = 45
$data[43,5] = 58
$data[43,6] = 103
$data[43,7] = 0.01205947781290247
$data[43,8] = 0.6645591851071303
$data[44,0] = 83
$data[44,1] = "Cromossoma 83"
$data[44,2] = 1410.088754026931
$data[44,3] = 15
$data[44,4] = 44
$data[44,5] = 58
$data[44,6] = 102
$data[44,7] = 0.01194239550403934
$data[44,8] = 0.6765015806111696
$data[45,0] = 25
$data[45,1] = "Cromossoma 25"
$data[45,2] = 1414.212193883911
$data[45,3] = 15
$data[45,4] = 42
$data[45,5] = 58
$data[45,6] = 100
$data[45,7] = 0.01170823088631308
$data[45,8] = 0.6882098114974827
$data[46,0] = 43
$data[46,1] = "Cromossoma 43"
$data[46,2] = 1424.760422311531
$data[46,3] = 15
$data[46,4] = 40
$data[46,5] = 58
$data[46,6] = 98
$data[46,7] = 0.01147406626858682
$data[46,8] = 0.6996838777660696
$data[47,0] = 15
$data[47,1] = "Cromossoma 15"
$data[47,2] = 1484.259807153725
$data[47,3] = 14
$data[47,4] = 14
$data[47,5] = 83
$data[47,6] = 97
$data[47,7] = 0.01135698395972369
$data[47,8] = 0.7110408617257933
$data[48,0] = 16
$data[48,1] = "Cromossoma 16"
$data[48,2] = 1425.404260637257
$data[48,3] = 15
$data[48,4] = 39
$data[48,5] = 58
$data[48,6] = 97
$data[48,7] = 0.01135698395972369
$data[48,8] = 0.722397845685517
$data[49,0] = 14
$data[49,1] = "Cromossoma 14"
$data[49,2] = 1428.180999308006
$data[49,3] = 15
$data[49,4] = 38
$data[49,5] = 58
$data[49,6] = 96
$data[49,7] = 0.01123990165086055
$data[49,8] = 0.7336377473363775
$data[50,0] = 3
$data[50,1] = "Cromossoma 3"
$data[50,2] = 1488.503263784041
$data[50,3] = 14
$data[50,4] = 12
$data[50,5] = 83
$data[50,6] = 95
$data[50,7] = 0.01112281934199742
$data[50,8] = 0.744760566678375
$data[51,0] = 54
$data[51,1] = "Cromossoma 54"
$data[51,2] = 1270.79645268163
$data[51,3] = 16
$data[51,4] = 81
$data[51,5] = 12
$data[51,6] = 93
$data[51,7] = 0.01088865472427116
$data[51,8] = 0.7556492214026461
$data[52,0] = 5
$data[52,1] = "Cromossoma 5"
$data[52,2] = 1436.454405949466
$data[52,3] = 15
$data[52,4] = 35
$data[52,5] = 58
$data[52,6] = 93
$data[52,7] = 0.01088865472427116
$data[52,8] = 0.7665378761269173
$data[53,0] = 21
$data[53,1] = "Cromossoma 21"
$data[53,2] = 1437.615342962999
$data[53,3] = 15
$data[53,4] = 34
$data[53,5] = 58
$data[53,6] = 92
$data[53,7] = 0.01077157241540803
$data[53,8] = 0.7773094485423253
$data[54,0] = 39
$data[54,1] = "Cromossoma 39"
$data[54,2] = 1438.899952071488
$data[54,3] = 15
$data[54,4] = 33
$data[54,5] = 58
$data[54,6] = 91
$data[54,7] = 0.0106544901065449
$data[54,8] = 0.7879639386488702
$data[55,0] = 50
$data[55,1] = "Cromossoma 50"
$data[55,2] = 1512.791900001893
$data[55,3] = 14
$data[55,4] = 7
$data[55,5] = 83
$data[55,6] = 90
$data[55,7] = 0.01053740779768177
$data[55,8] = 0.7985013464465519
$data[56,0] = 38
$data[56,1] = "Cromossoma 38"
$data[56,2] = 1442.015640197786
$data[56,3] = 15
$data[56,4] = 32
$data[56,5] = 58
$data[56,6] = 90
$data[56,7] = 0.01053740779768177
$data[56,8] = 0.8090387542442337
$data[57,0] = 62
$data[57,1] = "Cromossoma 62"
$data[57,2] = 1513.91418448466
$data[57,3] = 14
$data[57,4] = 6
$data[57,5] = 83
$data[57,6] = 89
$data[57,7] = 0.01042032548881864
$data[57,8] = 0.8194590797330523
$data[58,0] = 27
$data[58,1] = "Cromossoma 27"
$data[58,2] = 1445.457348595379
$data[58,3] = 15
$data[58,4] = 31
$data[58,5] = 58
$data[58,6] = 89
$data[58,7] = 0.01042032548881864
$data[58,8] = 0.829879405221871
$data[59,0] = 47
$data[59,1] = "Cromossoma 47"
$data[59,2] = 1447.578635610252
$data[59,3] = 15
$data[59,4] = 30
$data[59,5] = 58
$data[59,6] = 88
$data[59,7] = 0.01030324317995551
$data[59,8] = 0.8401826484018264
$data[60,0] = 34
$data[60,1] = "Cromossoma 34"
$data[60,2] = 1449.060109721396
$data[60,3] = 15
$data[60,4] = 29
$data[60,5] = 58
$data[60,6] = 87
$data[60,7] = 0.01018616087109238
$data[60,8] = 0.8503688092729188
$data[61,0] = 67
$data[61,1] = "Cromossoma 67"
$data[61,2] = 1449.509449272464
$data[61,3] = 15
$data[61,4] = 28
$data[61,5] = 58
$data[61,6] = 86
$data[61,7] = 0.01006907856222925
$data[61,8] = 0.8604378878351481
$data[62,0] = 52
$data[62,1] = "Cromossoma 52"
$data[62,2] = 1350.678935605421
$data[62,3] = 16
$data[62,4] = 73
$data[62,5] = 12
$data[62,6] = 85
$data[62,7] = 0.009951996253366116
$data[62,8] = 0.8703898840885141
$data[63,0] = 7
$data[63,1] = "Cromossoma 7"
$data[63,2] = 1452.20371323302
$data[63,3] = 15
$data[63,4] = 25
$data[63,5] = 58
$data[63,6] = 83
$data[63,7] = 0.009717831635639855
$data[63,8] = 0.880107715724154
$data[64,0] = 30
$data[64,1] = "Cromossoma 30"
$data[64,2] = 1463.001962675017
$data[64,3] = 15
$data[64,4] = 23
$data[64,5] = 58
$data[64,6] = 81
$data[64,7] = 0.009483667017913594
$data[64,8] = 0.8895913827420676
$data[65,0] = 31
$data[65,1] = "Cromossoma 31"
$data[65,2] = 1475.638097178263
$data[65,3] = 15
$data[65,4] = 19
$data[65,5] = 58
$data[65,6] = 77
$data[65,7] = 0.00901533778246107
$data[65,8] = 0.8986067205245287
$data[66,0] = 26
$data[66,1] = "Cromossoma 26"
$data[66,2] = 1477.491985185871
$data[66,3] = 15
$data[66,4] = 17
$data[66,5] = 58
$data[66,6] = 75
$data[66,7] = 0.008781173164734809
$data[66,8] = 0.9073878936892635
$data[67,0] = 48
$data[67,1] = "Cromossoma 48"
$data[67,2] = 1378.9851363771
$data[67,3] = 16
$data[67,4] = 63
$data[67,5] = 12
$data[67,6] = 75
$data[67,7] = 0.008781173164734809
$data[67,8] = 0.9161690668539983
$data[68,0] = 4
$data[68,1] = "Cromossoma 4"
$data[68,2] = 1480.022800173836
$data[68,3] = 15
$data[68,4] = 16
$data[68,5] = 58
$data[68,6] = 74
$data[68,7] = 0.008664090855871678
$data[68,8] = 0.9248331577098701
$data[69,0] = 32
$data[69,1] = "Cromossoma 32"
$data[69,2] = 1482.060005881544
$data[69,3] = 15
$data[69,4] = 15
$data[69,5] = 58
$data[69,6] = 73
$data[69,7] = 0.008547008547008548
$data[69,8] = 0.9333801662568786
$data[70,0] = 81
$data[70,1] = "Cromossoma 81"
$data[70,2] = 1496.239189270172
$data[70,3] = 15
$data[70,4] = 10
$data[70,5] = 58
$data[70,6] = 68
$data[70,7] = 0.007961597002692894
$data[70,8] = 0.9413417632595715
$data[71,0] = 19
$data[71,1] = "Cromossoma 19"
$data[71,2] = 1498.905199353085
$data[71,3] = 15
$data[71,4] = 9
$data[71,5] = 58
$data[71,6] = 67
$data[71,7] = 0.007844514693829763
$data[71,8] = 0.9491862779534013
$data[72,0] = 68
$data[72,1] = "Cromossoma 68"
$data[72,2] = 1519.718394030406
$data[72,3] = 15
$data[72,4] = 4
$data[72,5] = 58
$data[72,6] = 62
$data[72,7] = 0.007259103149514109
$data[72,8] = 0.9564453811029153
$data[73,0] = 69
$data[73,1] = "Cromossoma 69"
$data[73,2] = 1521.06343330563
$data[73,3] = 15
$data[73,4] = 3
$data[73,5] = 58
$data[73,6] = 61
$data[73,7] = 0.007142020840650977
$data[73,8] = 0.9635874019435663
$data[74,0] = 66
$data[74,1] = "Cromossoma 66"
$data[74,2] = 1528.018775571936
$data[74,3] = 15
$data[74,4] = 2
$data[74,5] = 58
$data[74,6] = 60
$data[74,7] = 0.007024938531787847
$data[74,8] = 0.9706123404753542
$data[75,0] = 73
$data[75,1] = "Cromossoma 73"
$data[75,2] = 1415.397017225511
$data[75,3] = 16
$data[75,4] = 41
$data[75,5] = 12
$data[75,6] = 53
$data[75,7] = 0.006205362369745932
$data[75,8] = 0.9768177028451001
$data[76,0] = 80
$data[76,1] = "Cromossoma 80"
$data[76,2] = 1460.579789035702
$data[76,3] = 16
$data[76,4] = 24
$data[76,5] = 12
$data[76,6] = 36
$data[76,7] = 0.004214963119072708
$data[76,8] = 0.9810326659641727
$data[77,0] = 75
$data[77,1] = "Cromossoma 75"
$data[77,2] = 1469.41905046335
$data[77,3] = 16
$data[77,4] = 22
$data[77,5] = 12
$data[77,6] = 34
$data[77,7] = 0.003980798501346447
$data[77,8] = 0.9850134644655192
$data[78,0] = 37
$data[78,1] = "Cromossoma 37"
$data[78,2] = 1476.736592400702
$data[78,3] = 16
$data[78,4] = 18
$data[78,5] = 12
$data[78,6] = 30
$data[78,7] = 0.003512469265893923
$data[78,8] = 0.9885259337314131
$data[79,0] = 10
$data[79,1] = "Cromossoma 10"
$data[79,2] = 1486.222733802868
$data[79,3] = 16
$data[79,4] = 13
$data[79,5] = 12
$data[79,6] = 25
$data[79,7] = 0.00292705772157827
$data[79,8] = 0.9914529914529914
$data[80,0] = 1
$data[80,1] = "Cromossoma 1"
$data[80,2] = 1489.511405157817
$data[80,3] = 16
$data[80,4] = 11
$data[80,5] = 12
$data[80,6] = 23
$data[80,7] = 0.002692893103852008
$data[80,8] = 0.9941458845568434
$data[81,0] = 17
$data[81,1] = "Cromossoma 17"
$data[81,2] = 1504.517832042983
$data[81,3] = 16
$data[81,4] = 8
$data[81,5] = 12
$data[81,6] = 20
$data[81,7] = 0.002341646177262616
$data[81,8] = 0.996487530734106
$data[82,0] = 45
$data[82,1] = "Cromossoma 45"
$data[82,2] = 1514.568056234249
$data[82,3] = 16
$data[82,4] = 5
$data[82,5] = 12
$data[82,6] = 17
$data[82,7] = 0.001990399250673223
$data[82,8] = 0.9984779299847792
$data[83,0] = 12
$data[83,1] = "Cromossoma 12"
$data[83,2] = 1530.369689696431
$data[83,3] = 16
$data[83,4] = 1
$data[83,5] = 12
$data[83,6] = 13
$data[83,7] = 0.0015220700152207
$data[83,8] = 1
$ws.Range("A2:I85").Value = $data
